# aggiornamento fino al 26/03
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows to append at the bottom of the table (rows 234-238)
$newRows = @(
    @{ Row = 234; Date = 44308; B = 0; C = 1; D = 28.87669650591972 },
    @{ Row = 235; Date = 44309; B = 1; C = 2; D = 57.75339301183945 },
    @{ Row = 236; Date = 44310; B = 0; C = 2; D = 57.75339301183945 },
    @{ Row = 237; Date = 44311; B = 0; C = 2; D = 57.75339301183945 },
    @{ Row = 238; Date = 44312; B = 0; C = 2; D = 57.75339301183945 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $prevRow = $row - 1

    # Copy formatting (date style with border/bold/centered alignment) from the
    # column A cell directly above so the new cell matches the existing table style.
    $ws.Range("A$prevRow").Copy($ws.Range("A$row"))

    $ws.Range("A$row").Value = $r.Date
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
}
